# Refresh the cryptos table: update Price (col D) and Volume(1h) (col E)
# for each coin row (2-51) on the active worksheet, per the latest pull.
#
# D-column values are plain decimal-looking strings (e.g. "357.66") that
# Excel auto-coerces to numbers on assignment; the source data stores them
# as text (preserves formatting like trailing zeros / thousands dots), so
# each D write is wrapped with a Text number format, then the style is
# reset to "Normal" so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "51.767.09"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  -0.31%  "
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "2.771.49"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("E4").Value = "  +0.07%  "
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "357.66"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("E6").Value = "  -4.28%  "
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.96%  "
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "39.96"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  -4.57%  "
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0852"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").Value = "  +0.20%  "
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "19.44"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  -2.16%  "
$ws.Range("E14").Value = "  -2.32%  "
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "3.208.34"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  -1.70%  "
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "2.743.40"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  -3.11%  "
$ws.Range("E17").Value = "  +2.40%  "
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "51.714.63"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("E19").Value = "  +1.51%  "
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "3.11"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("E21").Value = "  -3.98%  "
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0980"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  -0.14%  "
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "273.18"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  +1.76%  "
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "69.57"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").Value = "  -1.69%  "
$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "26.46"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  -1.58%  "
$ws.Range("E27").Value = "  -0.01%  "
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "10.14"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  -1.08%  "
$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "2.19"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  -2.55%  "
$ws.Range("E30").Value = "  +0.90%  "
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0472"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  +7.13%  "
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "51.23"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  +0.67%  "
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "33.93"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  +0.11%  "
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "5.74"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  -1.24%  "
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "5.31"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  +7.74%  "
$ws.Range("E36").Value = "  +1.06%  "
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "1.00"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  +0.13%  "
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "3.19"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("E39").Value = "  -5.12%  "
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "18.02"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -2.09%  "
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "2.62"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  +3.26%  "
$ws.Range("E42").Value = "  -0.63%  "
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "124.97"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("E44").Value = "  -0.99%  "
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "21.89"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -6.31%  "
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "2.056.77"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  +0.07%  "
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "2.31"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  +0.10%  "
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "3.22"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  -3.95%  "
$ws.Range("E49").Value = "  +0.67%  "
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "0.940"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  +0.52%  "
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "8.98"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  +0.22%  "
